# Bugfix/Update wrong column names (#2)
#
# Fix the mis-typed "Provider _ID" header (stray space before the
# underscore) so it reads "Provider__ID", and restore the sheet's
# scroll position / selection back to the top-left of the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content fix --------------------------------------------------------
# T1 currently holds the shared string "Provider _ID" -> "Provider__ID"
$ws.Range("T1").Value = "Provider__ID"

# --- View / selection fix -----------------------------------------------
# Scroll the window back so column A is the left-most visible column
# (topLeftCell was "X1", should be "A1") and move the active selection
# to C11 (was AE1).
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("C11").Select()
